$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final 10 data rows (A2:T11) for the sCs/FAPs/Calca/Calcrl/ECs/M1/M2 LR-pairs table,
# following Dr Hou's revised NATMI computation advice.
$rows = @(
    @("FAPs","Calca","Calcrl","ECs",1,0.3333333333333333,0.06082199999999999,0.182466,0.1419252066466042,0.1419252066466042,3,1,66.996675,200.990025,0.6276213057730995,0.6276213057730995,4.074871766849999,36.67384590165,0.08907528351765868,0.08907528351765868),
    @("FAPs","Calca","Calcrl","FAPs",1,0.3333333333333333,0.06082199999999999,0.182466,0.1419252066466042,0.1419252066466042,3,1,13.59324033333333,40.779721,0.127340756055345,0.127340756055345,0.8267680635539999,7.440912571985999,0.01807286311768964,0.01807286311768964),
    @("FAPs","Calca","Calcrl","M1",1,0.3333333333333333,0.06082199999999999,0.182466,0.1419252066466042,0.1419252066466042,3,1,12.031188,36.093564,0.1127075324642849,0.1127075324642849,0.731760916536,6.585848248824,0.01599603983562248,0.01599603983562248),
    @("FAPs","Calca","Calcrl","M2",1,0.3333333333333333,0.06082199999999999,0.182466,0.1419252066466042,0.1419252066466042,3,1,12.55587433333333,37.667623,0.1176227662672754,0.1176227662672754,0.763673388702,6.873060498317999,0.01669363540882828,0.01669363540882828),
    @("FAPs","Calca","Calcrl","sCs",1,0.3333333333333333,0.06082199999999999,0.182466,0.1419252066466042,0.1419252066466042,3,1,1.569996,4.709988000000001,0.01470763943999524,0.01470763943999524,0.09549029671200002,0.8594126704080001,0.002087384766805071,0.00208738476680507),
    @("sCs","Calca","Calcrl","ECs",3,1,0.3677276666666667,1.103183,0.8580747933533959,0.8580747933533959,3,1,66.996675,200.990025,0.6276213057730995,0.6276213057730995,24.636530972175,221.728778749575,0.5385460222554408,0.5385460222554408),
    @("sCs","Calca","Calcrl","FAPs",3,1,0.3677276666666667,1.103183,0.8580747933533959,0.8580747933533959,3,1,13.59324033333333,40.779721,0.127340756055345,0.127340756055345,4.998610550215888,44.987494951943,0.1092678929376553,0.1092678929376553),
    @("sCs","Calca","Calcrl","M1",3,1,0.3677276666666667,1.103183,0.8580747933533959,0.8580747933533959,3,1,12.031188,36.093564,0.1127075324642849,0.1127075324642849,4.424200690468,39.817806214212,0.09671149262866244,0.09671149262866244),
    @("sCs","Calca","Calcrl","M2",3,1,0.3677276666666667,1.103183,0.8580747933533959,0.8580747933533959,3,1,12.55587433333333,37.667623,0.1176227662672754,0.1176227662672754,4.617142371556556,41.554281344009,0.1009291308584471,0.1009291308584471),
    @("sCs","Calca","Calcrl","sCs",3,1,0.3677276666666667,1.103183,0.8580747933533959,0.8580747933533959,3,1,1.569996,4.709988000000001,0.01470763943999524,0.01470763943999524,0.5773309657560002,5.195978691804001,0.01262025467319017,0.01262025467319017)
)

$nr = $rows.Count
$nc = $rows[0].Count
$arr = New-Object "object[,]" $nr,$nc
for ($i = 0; $i -lt $nr; $i++) {
    for ($j = 0; $j -lt $nc; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}

$ws.Range("A2:T11").Value = $arr
